$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.454.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "'2.599.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'520.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "'144.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("D9").Value = "'2.617.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("D10").Value = "'6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").Value = "'0.327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.22%  "

$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").Value = "'3.058.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").Value = "'58.343.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").Value = "'20.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "'2.608.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "

$ws.Range("D19").Value = "'341.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.20%  "

$ws.Range("D20").Value = "'4.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("E21").Value = "  +0.59%  "

$ws.Range("D22").Value = "'6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.56%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "'66.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("D25").Value = "'0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("D26").Value = "'0.404"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'2.718.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").Value = "'7.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("D30").Value = "'0.0₃0753"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.49%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "'6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "

$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").Value = "'18.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").Value = "'149.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("D36").Value = "'4.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("E37").Value = "  -1.29%  "

$ws.Range("D38").Value = "'0.871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("D39").Value = "'0.846"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "

$ws.Range("E40").Value = "  +2.40%  "

$ws.Range("D41").Value = "'36.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("D42").Value = "'3.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("D43").Value = "'278.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.10%  "

$ws.Range("D44").Value = "'0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").Value = "'0.599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.39%  "

$ws.Range("D46").Value = "'0.0955"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").Value = "'10.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").Value = "'18.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("D49").Value = "'0.0524"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'19.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.26%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.987.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.95%  "
